{"js": "// Replace the division-problem text runs in the worksheet table with the\n// newly generated problems. Each \"before\" expression is unique within the\n// document, so an exact, case-sensitive search-and-replace on the full text\n// of each cell's run is sufficient and safe.\nconst replacements = [\n  [\"90\u00f73=\", \"24\u00f77=\"],\n  [\"85\u00f76=\", \"82\u00f78=\"],\n  [\"96\u00f77=\", \"27\u00f78=\"],\n  [\"22\u00f73=\", \"79\u00f75=\"],\n  [\"85\u00f72=\", \"54\u00f77=\"],\n  [\"21\u00f72=\", \"23\u00f76=\"],\n  [\"36\u00f74=\", \"45\u00f73=\"],\n  [\"88\u00f75=\", \"43\u00f76=\"],\n  [\"65\u00f73=\", \"43\u00f74=\"],\n  [\"29\u00f74=\", \"25\u00f76=\"],\n  [\"80\u00f74=\", \"51\u00f73=\"],\n  [\"16\u00f72=\", \"11\u00f74=\"],\n  [\"30\u00f76=\", \"10\u00f72=\"],\n  [\"48\u00f72=\", \"68\u00f72=\"],\n  [\"82\u00f73=\", \"52\u00f76=\"],\n  [\"57\u00f74=\", \"52\u00f75=\"],\n  [\"19\u00f78=\", \"49\u00f76=\"],\n  [\"71\u00f73=\", \"88\u00f75=\"],\n  [\"71\u00f77=\", \"20\u00f74=\"],\n  [\"93\u00f75=\", \"36\u00f75=\"],\n  [\"87\u00f79=\", \"24\u00f74=\"],\n  [\"36\u00f72=\", \"75\u00f79=\"],\n  [\"64\u00f79=\", \"28\u00f72=\"],\n  [\"66\u00f76=\", \"24\u00f73=\"],\n  [\"59\u00f73=\", \"95\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text runs in the worksheet table with the\n# newly generated problems. Each \"before\" expression is unique within the\n# document, so an exact, case-sensitive Find/Replace on each one is safe and\n# will touch exactly the single matching cell.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{Before=\"90\u00f73=\"; After=\"24\u00f77=\"},\n  @{Before=\"85\u00f76=\"; After=\"82\u00f78=\"},\n  @{Before=\"96\u00f77=\"; After=\"27\u00f78=\"},\n  @{Before=\"22\u00f73=\"; After=\"79\u00f75=\"},\n  @{Before=\"85\u00f72=\"; After=\"54\u00f77=\"},\n  @{Before=\"21\u00f72=\"; After=\"23\u00f76=\"},\n  @{Before=\"36\u00f74=\"; After=\"45\u00f73=\"},\n  @{Before=\"88\u00f75=\"; After=\"43\u00f76=\"},\n  @{Before=\"65\u00f73=\"; After=\"43\u00f74=\"},\n  @{Before=\"29\u00f74=\"; After=\"25\u00f76=\"},\n  @{Before=\"80\u00f74=\"; After=\"51\u00f73=\"},\n  @{Before=\"16\u00f72=\"; After=\"11\u00f74=\"},\n  @{Before=\"30\u00f76=\"; After=\"10\u00f72=\"},\n  @{Before=\"48\u00f72=\"; After=\"68\u00f72=\"},\n  @{Before=\"82\u00f73=\"; After=\"52\u00f76=\"},\n  @{Before=\"57\u00f74=\"; After=\"52\u00f75=\"},\n  @{Before=\"19\u00f78=\"; After=\"49\u00f76=\"},\n  @{Before=\"71\u00f73=\"; After=\"88\u00f75=\"},\n  @{Before=\"71\u00f77=\"; After=\"20\u00f74=\"},\n  @{Before=\"93\u00f75=\"; After=\"36\u00f75=\"},\n  @{Before=\"87\u00f79=\"; After=\"24\u00f74=\"},\n  @{Before=\"36\u00f72=\"; After=\"75\u00f79=\"},\n  @{Before=\"64\u00f79=\"; After=\"28\u00f72=\"},\n  @{Before=\"66\u00f76=\"; After=\"24\u00f73=\"},\n  @{Before=\"59\u00f73=\"; After=\"95\u00f77=\"}\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair.Before\n  $find.Replacement.Text = $pair.After\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($pair.Before, $false, $true, $false, $false, $false, $true, 1, $false, $pair.After, 2)\n}\n"}
